$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 10, shifting existing rows 10-21 down to 12-23
$ws.Rows("10:11").Insert()

# Populate new row 10 (Primera) - duplicate of old row 10 with updated date
$ws.Cells.Item(10,1).Value = 11
$ws.Cells.Item(10,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(10,3).Value = "Bíobío"
$ws.Cells.Item(10,4).Value = 44554
$ws.Cells.Item(10,5).Value = 8
$ws.Cells.Item(10,6).Value = 100112037
$ws.Cells.Item(10,7).Value = "Cebollín"
$ws.Cells.Item(10,8).Value = "Sin especificar"
$ws.Cells.Item(10,9).Value = "Primera"
$ws.Cells.Item(10,10).Value = 200
$ws.Cells.Item(10,11).Value = 600
$ws.Cells.Item(10,12).Value = 700
$ws.Cells.Item(10,13).Value = 650
$ws.Cells.Item(10,14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(10,15).Value = "Región de Ñuble"
$ws.Cells.Item(10,16).Value = 108
$ws.Cells.Item(10,17).Value = 6
$ws.Cells.Item(10,18).Value = "Hortaliza"

# Populate new row 11 (Segunda) - duplicate of old row 11 with updated date
$ws.Cells.Item(11,1).Value = 11
$ws.Cells.Item(11,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(11,3).Value = "Bíobío"
$ws.Cells.Item(11,4).Value = 44554
$ws.Cells.Item(11,5).Value = 8
$ws.Cells.Item(11,6).Value = 100112037
$ws.Cells.Item(11,7).Value = "Cebollín"
$ws.Cells.Item(11,8).Value = "Sin especificar"
$ws.Cells.Item(11,9).Value = "Segunda"
$ws.Cells.Item(11,10).Value = 100
$ws.Cells.Item(11,11).Value = 500
$ws.Cells.Item(11,12).Value = 500
$ws.Cells.Item(11,13).Value = 500
$ws.Cells.Item(11,14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(11,15).Value = "Región de Ñuble"
$ws.Cells.Item(11,16).Value = 83
$ws.Cells.Item(11,17).Value = 6
$ws.Cells.Item(11,18).Value = "Hortaliza"
